$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 rework -------------------------------------------------------
# A3 used to hold formula "=1+0" -> keep the value (1) but drop the formula.
$ws.Range("A3").Value = 1

# B3 used to hold formula CONCATENATE("a","b") -> that content moves down
# to B4 (new row). Clear B3 entirely for now; it gets re-created below.
$ws.Range("B3").ClearContents()

# C3 used to hold formula ("a"="b") evaluating to FALSE -> now a literal
# TRUE boolean value (no formula).
$ws.Range("C3").Value = $true

# D3 keeps its existing formula (=D2) untouched.

# --- New row 4 -----------------------------------------------------------
$ws.Range("A4").Value = 1
$ws.Range("B4").Formula = '=CONCATENATE("a", "b")'

# --- New column E: per-row type marker ----------------------------------
$ws.Range("E1").Value = "string_in_row_3"
$ws.Range("E2").Value = 1
$ws.Range("E3").Value = 2
$ws.Range("E4").Value = "ab"

# --- Selection moves to E5 (just past the used range) -------------------
$ws.Range("E5").Select()
